$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.583.74'
$ws.Range("E2").Value = '  +4.95%  '
$ws.Range("D3").Value = '1.592.60'
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.69%  '
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.13'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +9.51%  '
$ws.Range("E9").Value = '  +1.14%  '
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("D12").Value = '1.819.68'
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").Value = '1.609.78'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.80'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.28%  '
$ws.Range("E15").Value = '  +2.42%  '
$ws.Range("D16").Value = '28.550.41'
$ws.Range("E16").Value = '  +5.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.05'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '233.33'
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.53'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.45'
$ws.Range("D23").ClearFormats()
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.70'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.34'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.64'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  +1.02%  '
$ws.Range("E29").Value = '  -0.66%  '
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.26'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.15'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").Value = '1.419.37'
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("E36").Value = '  -5.98%  '
$ws.Range("E37").Value = '  -0.57%  '
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("E39").Value = '  +8.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.544'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.816'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.73'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.34%  '
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.979'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.27%  '
$ws.Range("E45").Value = '  +6.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.79'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("D47").Value = '1.732.47'
$ws.Range("E47").Value = '  +1.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.01'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +2.31%  '
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0524'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '40.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +17.79%  '
